$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.896501666666667
$ws.Range("H2").Value = 11.689505
$ws.Range("I2").Value = 0.401720501899026
$ws.Range("J2").Value = 0.401720501899026
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 17.00372766666667
$ws.Range("N2").Value = 51.011183
$ws.Range("O2").Value = 0.6620200065567142
$ws.Range("P2").Value = 0.6620200065567141
$ws.Range("Q2").Value = 66.25505319271278
$ws.Range("R2").Value = 596.2954787344149
$ws.Range("S2").Value = 0.2659470093011597
$ws.Range("T2").Value = 0.2659470093011597
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.896501666666667
$ws.Range("H3").Value = 11.689505
$ws.Range("I3").Value = 0.401720501899026
$ws.Range("J3").Value = 0.401720501899026
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.393811333333333
$ws.Range("N3").Value = 16.181434
$ws.Range("O3").Value = 0.2100016587103467
$ws.Range("P3").Value = 0.2100016587103467
$ws.Range("Q3").Value = 21.01699485001889
$ws.Range("R3").Value = 189.15295365017
$ws.Range("S3").Value = 0.08436197173674842
$ws.Range("T3").Value = 0.08436197173674842
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.896501666666667
$ws.Range("H4").Value = 11.689505
$ws.Range("I4").Value = 0.401720501899026
$ws.Range("J4").Value = 0.401720501899026
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.287074
$ws.Range("N4").Value = 9.861222
$ws.Range("O4").Value = 0.1279783347329392
$ws.Range("P4").Value = 0.1279783347329391
$ws.Range("Q4").Value = 12.80808931945667
$ws.Range("R4").Value = 115.27280387511
$ws.Range("S4").Value = 0.05141152086111787
$ws.Range("T4").Value = 0.05141152086111785
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.23724
$ws.Range("H5").Value = 9.71172
$ws.Range("I5").Value = 0.3337521163387849
$ws.Range("J5").Value = 0.3337521163387849
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.00372766666667
$ws.Range("N5").Value = 51.011183
$ws.Range("O5").Value = 0.6620200065567142
$ws.Range("P5").Value = 0.6620200065567141
$ws.Range("Q5").Value = 55.04514735164
$ws.Range("R5").Value = 495.4063261647599
$ws.Range("S5").Value = 0.2209505782469197
$ws.Range("T5").Value = 0.2209505782469196
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.23724
$ws.Range("H6").Value = 9.71172
$ws.Range("I6").Value = 0.3337521163387849
$ws.Range("J6").Value = 0.3337521163387849
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.393811333333333
$ws.Range("N6").Value = 16.181434
$ws.Range("O6").Value = 0.2100016587103467
$ws.Range("P6").Value = 0.2100016587103467
$ws.Range("Q6").Value = 17.46106180072
$ws.Range("R6").Value = 157.14955620648
$ws.Range("S6").Value = 0.07008849802923342
$ws.Range("T6").Value = 0.07008849802923343
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.23724
$ws.Range("H7").Value = 9.71172
$ws.Range("I7").Value = 0.3337521163387849
$ws.Range("J7").Value = 0.3337521163387849
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.287074
$ws.Range("N7").Value = 9.861222
$ws.Range("O7").Value = 0.1279783347329392
$ws.Range("P7").Value = 0.1279783347329391
$ws.Range("Q7").Value = 10.64104743576
$ws.Range("R7").Value = 95.76942692183999
$ws.Range("S7").Value = 0.04271304006263187
$ws.Range("T7").Value = 0.04271304006263186
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.565792333333333
$ws.Range("H8").Value = 7.697377
$ws.Range("I8").Value = 0.2645273817621892
$ws.Range("J8").Value = 0.2645273817621892
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 17.00372766666667
$ws.Range("N8").Value = 51.011183
$ws.Range("O8").Value = 0.6620200065567142
$ws.Range("P8").Value = 0.6620200065567141
$ws.Range("Q8").Value = 43.62803408522122
$ws.Range("R8").Value = 392.652306766991
$ws.Range("S8").Value = 0.1751224190086349
$ws.Range("T8").Value = 0.1751224190086349
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.565792333333333
$ws.Range("H9").Value = 7.697377
$ws.Range("I9").Value = 0.2645273817621892
$ws.Range("J9").Value = 0.2645273817621892
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.393811333333333
$ws.Range("N9").Value = 16.181434
$ws.Range("O9").Value = 0.2100016587103467
$ws.Range("P9").Value = 0.2100016587103467
$ws.Range("Q9").Value = 13.83939976651311
$ws.Range("R9").Value = 124.554597898618
$ws.Range("S9").Value = 0.05555118894436482
$ws.Range("T9").Value = 0.05555118894436483
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.565792333333333
$ws.Range("H10").Value = 7.697377
$ws.Range("I10").Value = 0.2645273817621892
$ws.Range("J10").Value = 0.2645273817621892
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.287074
$ws.Range("N10").Value = 9.861222
$ws.Range("O10").Value = 0.1279783347329392
$ws.Range("P10").Value = 0.1279783347329391
$ws.Range("Q10").Value = 8.433949268299333
$ws.Range("R10").Value = 75.90554341469399
$ws.Range("S10").Value = 0.03385377380918943
$ws.Range("T10").Value = 0.03385377380918943
